# Fill in the new translation rows (4-15) that were pasted into the table,
# matching the "Шеймин 2" sheet upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 90
$ws.Range("C4").Value = ' I heard a [CS:I]Gracidea[CR] was used to\nget down from the summit.'
$ws.Range("D4").Value = ' Я узнал, что [CS:I]Грацидеей[CR]\nвоспользовались, чтобы спуститься с\nвершины горы.'
$ws.Range("E4").Value = ' Ÿ ôèîàì, œóï [CS:I]Ãñàøéäååê[CR]\nâïòðïìûèïâàìéòû, œóïáú òðôòóéóûòÿ ò\nâåñšéîú ãïñú.'
$ws.Rows.Item(4).RowHeight = 31.8

# Row 5
$ws.Range("B5").Value = 93
$ws.Range("C5").Value = ' That\''s right. Using this flower,\nwe can change our appearance.'
$ws.Range("D5").Value = ' Так и есть. Коснувшись этого\nцветка, мы можем изменить свой облик.'
$ws.Range("E5").Value = ' Óàë é åòóû. Ëïòîôâšéòû üóïãï\nøâåóëà, íú íïçåí éèíåîéóû òâïê ïáìéë.'
$ws.Rows.Item(5).RowHeight = 31.8

# Row 6
$ws.Range("B6").Value = 96
$ws.Range("C6").Value = ' Our appearance changes pretty\ndrastically, so please don\''t be surprised.'
$ws.Range("D6").Value = ' Наши облики сильно изменяются,\nпоэтому, не стоит этому удивляться.'
$ws.Range("E6").Value = ' Îàšé ïáìéëé òéìûîï éèíåîÿýóòÿ,\nðïüóïíô, îå òóïéó üóïíô ôäéâìÿóûòÿ.'
$ws.Range("A6:B6").WrapText = $true
$ws.Range("A6:B6").Borders.Item(9).LineStyle = 1
$ws.Range("A6:B6").Borders.Item(9).Weight = 2
$ws.Range("C6:E6").WrapText = $true
$ws.Range("C6:E6").Font.Size = 8
$ws.Range("C6:E6").Borders.Item(9).LineStyle = 1
$ws.Range("C6:E6").Borders.Item(9).Weight = 2
$ws.Rows.Item(6).RowHeight = 31.8

# Row 7
$ws.Range("A7").Value = 'SCRIPT/D73P11A/us0202.ssb'
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = ' Would you like a [CS:I]Gracidea[CR]?'
$ws.Range("D7").Value = ' Хотите взять [CS:I]Грацидею[CR]?'
$ws.Range("E7").Value = ' Öïóéóå âèÿóû [CS:I]Ãñàøéäåý[CR]?'
$ws.Rows.Item(7).RowHeight = 43.2

# Row 8
$ws.Range("A8").Value = 'SCRIPT/D73P11A/us0301.ssb'
$ws.Range("B8").Value = 34
$ws.Range("C8").Value = 'Yes, please.'
$ws.Range("D8").Value = 'Да, пожалуйста.'
$ws.Range("E8").Value = 'Äà, ðïçàìôêòóà.'
$ws.Rows.Item(8).RowHeight = 43.2

# Row 9
$ws.Range("A9").Value = 'SCRIPT/D73P11A/us0402.ssb'
$ws.Range("B9").Value = 39
$ws.Range("C9").Value = ' Here you go.'
$ws.Range("D9").Value = ' Держите.'
$ws.Range("E9").Value = ' Äåñçéóå.'
$ws.Rows.Item(9).RowHeight = 43.2

# Row 10
$ws.Range("A10").Value = 'SCRIPT/D73P11A/us2001.ssb'
$ws.Range("B10").Value = 49
$ws.Range("C10").Value = '[CN][hero] received\n[CN]a [CS:I]Gracidea[CR]!'
$ws.Range("D10").Value = '[CN][hero] получает [CS:I]Грацидею[CR]!'
$ws.Range("E10").Value = '[CN][hero] ðïìôœàåó [CS:I]Ãñàøéäåý[CR]!'
$ws.Rows.Item(10).RowHeight = 43.2

# Row 11
$ws.Range("A11").Value = 'SCRIPT/D73P11A/us2101.ssb'
$ws.Range("B11").Value = 58
$ws.Range("C11").Value = ' Oh? It seems you\''re carrying too\nmany items.'
$ws.Range("D11").Value = ' Что? Кажется, у вас много вещей.'
$ws.Range("E11").Value = ' Œóï? Ëàçåóòÿ, ô âàò íîïãï âåþåê.'
$ws.Rows.Item(11).RowHeight = 43.2

# Row 12
$ws.Range("B12").Value = 65
$ws.Range("C12").Value = 'No, thanks.'
$ws.Range("D12").Value = 'Спасибо, не нужно.'
$ws.Range("E12").Value = 'Òðàòéáï, îå îôçîï.'

# Row 13
$ws.Range("B13").Value = 69
$ws.Range("C13").Value = ' I see.\nWell, you\''re still welcome to one if you change\nyour mind.'
$ws.Range("D13").Value = ' Понятно. Ну, если передумаете,\nобращайтесь.'
$ws.Range("E13").Value = ' Ðïîÿóîï. Îô, åòìé ðåñåäôíàåóå,\nïáñàþàêóåòû.'
$ws.Range("A13:B13").WrapText = $true
$ws.Range("A13:B13").Borders.Item(9).LineStyle = 1
$ws.Range("A13:B13").Borders.Item(9).Weight = 2
$ws.Range("C13:E13").WrapText = $true
$ws.Range("C13:E13").Font.Size = 8
$ws.Range("C13:E13").Borders.Item(9).LineStyle = 1
$ws.Range("C13:E13").Borders.Item(9).Weight = 2
$ws.Rows.Item(13).RowHeight = 21.6

# Row 14
$ws.Range("B14").Value = 19
$ws.Range("C14").Value = ' When we [CS:K]Shaymin[CR] use a\n[CS:I]Gracidea[CR], we can change our appearance.'
$ws.Range("D14").Value = ' Когда мы, [CS:K]Шеймины[CR], касаемся\n[CS:I]Грацидеи[CR], мы можем изменить нашу\nвнешность.'
$ws.Range("E14").Value = ' Ëïãäà íú, [CS:K]Šåêíéîú[CR], ëàòàåíòÿ\n[CS:I]Ãñàøéäåé[CR], íú íïçåí éèíåîéóû îàšô\nâîåšîïòóû.'
$ws.Rows.Item(14).RowHeight = 31.8

# Row 15
$ws.Range("B15").Value = 22
$ws.Range("C15").Value = ' Our appearance changes pretty\ndrastically, so please don\''t be surprised.'
$ws.Range("D15").Value = ' Наши облики сильно изменяются,\nпоэтому, не стоит этому удивляться.'
$ws.Range("E15").Value = ' Îàšé ïáìéëé òéìûîï éèíåîÿýóòÿ,\nðïüóïíô, îå òóïéó üóïíô ôäéâìÿóûòÿ.'
$ws.Range("A15:B15").WrapText = $true
$ws.Range("A15:B15").Borders.Item(9).LineStyle = 1
$ws.Range("A15:B15").Borders.Item(9).Weight = 2
$ws.Range("C15:E15").WrapText = $true
$ws.Range("C15:E15").Font.Size = 8
$ws.Range("C15:E15").Borders.Item(9).LineStyle = 1
$ws.Range("C15:E15").Borders.Item(9).Weight = 2
$ws.Rows.Item(15).RowHeight = 31.8

# Match the author's final scroll position / selection in the sheet view.
$ws.Range("D10").Select()
